# fix(data): handling missing data at the beginning and at the end of the series
#
# 1) Row 12 has a stray cell F12 (= "k") left over from a bad paste; it
#    doesn't belong in the series at all, so remove it outright. Clear()
#    (contents + formatting) drops the cell entirely instead of leaving an
#    empty placeholder behind, and - unlike Delete() - it does not shift
#    any neighbouring cells/rows/columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F12").Clear()

# 2) Column C (MSCI_ACWI) is missing for the last stretch of rows
#    (288-321) even though the other columns (B, D) have data there -
#    restore the missing MSCI_ACWI values.
$acwiValues = @{
    288 = 1300.54199
    289 = 1402.013192
    290 = 1347.399869
    291 = 1444.315908
    292 = 1403.429303
    293 = 1447.674908
    294 = 1469.169928
    295 = 1454.481784
    296 = 1539.573475
    297 = 1596.391861
    298 = 1552.416902
    299 = 1488.773934
    300 = 1444.344327
    301 = 1578.338665
    302 = 1654.705023
    303 = 1664.763092
    304 = 1736.815265
    305 = 1792.3067
    306 = 1733.93893
    307 = 1805.457166
    308 = 1846.320461
    309 = 1876.606339
    310 = 1924.907114
    311 = 1970.352877
    312 = 1926.776187
    313 = 1999.474831
    314 = 1952.803514
    315 = 2018.786618
    316 = 2007.191047
    317 = 1928.980371
    318 = 1947.837353
    319 = 2061.073367
    320 = 2154.502262
    321 = 2184.265977
}

foreach ($row in $acwiValues.Keys) {
    $ws.Cells.Item([int]$row, 3).Value = $acwiValues[$row]
}

# 3) Update the saved view state: scroll back to the top and move the
#    selection to G13.
$ws.Activate()
$ws.Range("G13").Select() | Out-Null
